# Add thêm nhân sự Nguyễn Hữu Quang
# Updates the "Tháng 8" (row 6) financial figures and refreshes the
# last_edited_time stamps (column D) that Notion bumped when the row
# was edited to add the new hire.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tháng 8 (row 6) figures bumped up ---
$ws.Range("T6").Value = 143500000
$ws.Range("W6").Value = 253041000
$ws.Range("AA6").Value = 379659000
$ws.Range("AE6").Value = 632700000
$ws.Range("AH6").Value = 535500000
$ws.Range("AK6").Value = 81
$ws.Range("AQ6").Value = 679000000

# --- last_edited_time stamps refreshed ---
# Rows 6-9 (Tháng 8, 7, 6, 5) and rows 10-12 (Tháng 4, 3, 2) now share the
# newer timestamp; row 13 (Tháng 1) keeps the other refreshed timestamp.
$ws.Range("D6").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D7").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D8").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D9").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D10").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D11").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D12").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D13").Value = "2024-08-31T15:45:00.000Z"
